# Weekly update: insert the newest "Perejil" (Terminal Hortofrutícola Agro
# Chillán) price record at the top of the data block (row 72), pushing all
# existing records (old rows 72-99) down by one row (to 73-100).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 72; this shifts rows 72:99 down to 73:100
# and grows the sheet dimension from A1:R99 to A1:R100 automatically.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the latest week's record.
$ws.Range("A72").Value = 7
$ws.Range("B72").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C72").Value = "Ñuble"
$ws.Range("D72").Value = 45119
$ws.Range("E72").Value = 16
$ws.Range("F72").Value = 100112044
$ws.Range("G72").Value = "Perejil"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 200
$ws.Range("K72").Value = 1500
$ws.Range("L72").Value = 1500
$ws.Range("M72").Value = 1500
$ws.Range("N72").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O72").Value = "Región de Ñuble"
$ws.Range("P72").Value = 1500
$ws.Range("Q72").Value = 1
$ws.Range("R72").Value = "Hortaliza"
